$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F values for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1173
$ws1.Range("F4").Value = 274
$ws1.Range("F7").Value = 12298
$ws1.Range("F8").Value = 59
$ws1.Range("F9").Value = 13
$ws1.Range("F11").Value = 147
$ws1.Range("F12").Value = 12089
$ws1.Range("F13").Value = 4815
$ws1.Range("F14").Value = 4670
$ws1.Range("F15").Value = 124
$ws1.Range("F20").Value = 2

# Sheet "全部类型" (sheet4): update column F values for the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1173
$ws4.Range("F4").Value = 274
$ws4.Range("F9").Value = 12298
$ws4.Range("F10").Value = 59
$ws4.Range("F11").Value = 13
$ws4.Range("F13").Value = 147
$ws4.Range("F14").Value = 12089
$ws4.Range("F15").Value = 4815
$ws4.Range("F16").Value = 4670
$ws4.Range("F17").Value = 124
$ws4.Range("F22").Value = 2
